$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.135749
$ws.Range("H2").Value = 0.407247
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 0.8206224791043335
$ws.Range("R2").Value = 7.385602311939002
$ws.Range("S2").Value = 0.8160840232643366
$ws.Range("T2").Value = 0.8160840232643367

# Row 3
$ws.Range("G3").Value = 0.135749
$ws.Range("H3").Value = 0.407247
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.09264100005666669
$ws.Range("R3").Value = 0.8337690005100001
$ws.Range("S3").Value = 0.09212864864242169
$ws.Range("T3").Value = 0.09212864864242169

# Row 4
$ws.Range("G4").Value = 0.135749
$ws.Range("H4").Value = 0.407247
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.092297781335
$ws.Range("R4").Value = 0.8306800320150001
$ws.Range("S4").Value = 0.09178732809324164
$ws.Range("T4").Value = 0.09178732809324165
